# update scripts wuth new tpm
# Refresh the NATMI ligand-receptor edge table (Pglyrp1 -> Trem1) with the
# newly computed TPM-based values, and add the 5 additional "MuSCs" /
# "Resolving-Mac" sending-cluster rows that come from the updated run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,"A").Value = "ECs"
$ws.Cells.Item(2,"B").Value = "Pglyrp1"
$ws.Cells.Item(2,"C").Value = "Trem1"
$ws.Cells.Item(2,"D").Value = "ECs"
$ws.Cells.Item(2,"E").Value = 3
$ws.Cells.Item(2,"F").Value = 1
$ws.Cells.Item(2,"G").Value = 1.164414
$ws.Cells.Item(2,"H").Value = 3.493242
$ws.Cells.Item(2,"I").Value = 0.2111338079173259
$ws.Cells.Item(2,"J").Value = 0.2149356714351638
$ws.Cells.Item(2,"K").Value = 1
$ws.Cells.Item(2,"L").Value = 0.3333333333333333
$ws.Cells.Item(2,"M").Value = 0.006995333333333333
$ws.Cells.Item(2,"N").Value = 0.020986
$ws.Cells.Item(2,"O").Value = 0.003071634989180631
$ws.Cells.Item(2,"P").Value = 0.003071634989180632
$ws.Cells.Item(2,"Q").Value = 0.008145464068000001
$ws.Cells.Item(2,"R").Value = 0.07330917661200001
$ws.Cells.Item(2,"S").Value = 0.0006485259917978007
$ws.Cells.Item(2,"T").Value = 0.0006602039288032812

# Row 3
$ws.Cells.Item(3,"A").Value = "ECs"
$ws.Cells.Item(3,"B").Value = "Pglyrp1"
$ws.Cells.Item(3,"C").Value = "Trem1"
$ws.Cells.Item(3,"D").Value = "Inflammatory-Mac"
$ws.Cells.Item(3,"E").Value = 3
$ws.Cells.Item(3,"F").Value = 1
$ws.Cells.Item(3,"G").Value = 1.164414
$ws.Cells.Item(3,"H").Value = 3.493242
$ws.Cells.Item(3,"I").Value = 0.2111338079173259
$ws.Cells.Item(3,"J").Value = 0.2149356714351638
$ws.Cells.Item(3,"K").Value = 3
$ws.Cells.Item(3,"L").Value = 1
$ws.Cells.Item(3,"M").Value = 1.994843666666667
$ws.Cells.Item(3,"N").Value = 5.984531
$ws.Cells.Item(3,"O").Value = 0.8759313262859122
$ws.Cells.Item(3,"P").Value = 0.8759313262859124
$ws.Cells.Item(3,"Q").Value = 2.322823893278001
$ws.Cells.Item(3,"R").Value = 20.905415039502
$ws.Cells.Item(3,"S").Value = 0.1849387163928183
$ws.Cells.Item(3,"T").Value = 0.1882688877463561

# Row 4
$ws.Cells.Item(4,"A").Value = "ECs"
$ws.Cells.Item(4,"B").Value = "Pglyrp1"
$ws.Cells.Item(4,"C").Value = "Trem1"
$ws.Cells.Item(4,"D").Value = "Resolving-Mac"
$ws.Cells.Item(4,"E").Value = 3
$ws.Cells.Item(4,"F").Value = 1
$ws.Cells.Item(4,"G").Value = 1.164414
$ws.Cells.Item(4,"H").Value = 3.493242
$ws.Cells.Item(4,"I").Value = 0.2111338079173259
$ws.Cells.Item(4,"J").Value = 0.2149356714351638
$ws.Cells.Item(4,"K").Value = 2
$ws.Cells.Item(4,"L").Value = 0.6666666666666666
$ws.Cells.Item(4,"M").Value = 0.2755583333333333
$ws.Cells.Item(4,"N").Value = 0.826675
$ws.Cells.Item(4,"O").Value = 0.120997038724907
$ws.Cells.Item(4,"P").Value = 0.120997038724907
$ws.Cells.Item(4,"Q").Value = 0.32086398115
$ws.Cells.Item(4,"R").Value = 2.88777583035
$ws.Cells.Item(4,"S").Value = 0.02554656553270976
$ws.Cells.Item(4,"T").Value = 0.02600657976000441

# Row 5
$ws.Cells.Item(5,"A").Value = "FAPs"
$ws.Cells.Item(5,"B").Value = "Pglyrp1"
$ws.Cells.Item(5,"C").Value = "Trem1"
$ws.Cells.Item(5,"D").Value = "ECs"
$ws.Cells.Item(5,"E").Value = 3
$ws.Cells.Item(5,"F").Value = 1
$ws.Cells.Item(5,"G").Value = 2.141562333333333
$ws.Cells.Item(5,"H").Value = 6.424687
$ws.Cells.Item(5,"I").Value = 0.3883122414613532
$ws.Cells.Item(5,"J").Value = 0.395304537763421
$ws.Cells.Item(5,"K").Value = 1
$ws.Cells.Item(5,"L").Value = 0.3333333333333333
$ws.Cells.Item(5,"M").Value = 0.006995333333333333
$ws.Cells.Item(5,"N").Value = 0.020986
$ws.Cells.Item(5,"O").Value = 0.003071634989180631
$ws.Cells.Item(5,"P").Value = 0.003071634989180632
$ws.Cells.Item(5,"Q").Value = 0.01498094237577778
$ws.Cells.Item(5,"R").Value = 0.134828481382
$ws.Cells.Item(5,"S").Value = 0.001192753467599851
$ws.Cells.Item(5,"T").Value = 0.001214231249576

# Row 6
$ws.Cells.Item(6,"A").Value = "FAPs"
$ws.Cells.Item(6,"B").Value = "Pglyrp1"
$ws.Cells.Item(6,"C").Value = "Trem1"
$ws.Cells.Item(6,"D").Value = "Inflammatory-Mac"
$ws.Cells.Item(6,"E").Value = 3
$ws.Cells.Item(6,"F").Value = 1
$ws.Cells.Item(6,"G").Value = 2.141562333333333
$ws.Cells.Item(6,"H").Value = 6.424687
$ws.Cells.Item(6,"I").Value = 0.3883122414613532
$ws.Cells.Item(6,"J").Value = 0.395304537763421
$ws.Cells.Item(6,"K").Value = 3
$ws.Cells.Item(6,"L").Value = 1
$ws.Cells.Item(6,"M").Value = 1.994843666666667
$ws.Cells.Item(6,"N").Value = 5.984531
$ws.Cells.Item(6,"O").Value = 0.8759313262859122
$ws.Cells.Item(6,"P").Value = 0.8759313262859124
$ws.Cells.Item(6,"Q").Value = 4.27208205742189
$ws.Cells.Item(6,"R").Value = 38.44873851679701
$ws.Cells.Item(6,"S").Value = 0.3401348566762985
$ws.Cells.Item(6,"T").Value = 0.3462596280499529

# Row 7
$ws.Cells.Item(7,"A").Value = "FAPs"
$ws.Cells.Item(7,"B").Value = "Pglyrp1"
$ws.Cells.Item(7,"C").Value = "Trem1"
$ws.Cells.Item(7,"D").Value = "Resolving-Mac"
$ws.Cells.Item(7,"E").Value = 3
$ws.Cells.Item(7,"F").Value = 1
$ws.Cells.Item(7,"G").Value = 2.141562333333333
$ws.Cells.Item(7,"H").Value = 6.424687
$ws.Cells.Item(7,"I").Value = 0.3883122414613532
$ws.Cells.Item(7,"J").Value = 0.395304537763421
$ws.Cells.Item(7,"K").Value = 2
$ws.Cells.Item(7,"L").Value = 0.6666666666666666
$ws.Cells.Item(7,"M").Value = 0.2755583333333333
$ws.Cells.Item(7,"N").Value = 0.826675
$ws.Cells.Item(7,"O").Value = 0.120997038724907
$ws.Cells.Item(7,"P").Value = 0.120997038724907
$ws.Cells.Item(7,"Q").Value = 0.5901253473027778
$ws.Cells.Item(7,"R").Value = 5.311128125725001
$ws.Cells.Item(7,"S").Value = 0.04698463131745479
$ws.Cells.Item(7,"T").Value = 0.04783067846389213

# Row 8
$ws.Cells.Item(8,"A").Value = "Inflammatory-Mac"
$ws.Cells.Item(8,"B").Value = "Pglyrp1"
$ws.Cells.Item(8,"C").Value = "Trem1"
$ws.Cells.Item(8,"D").Value = "ECs"
$ws.Cells.Item(8,"E").Value = 3
$ws.Cells.Item(8,"F").Value = 1
$ws.Cells.Item(8,"G").Value = 1.198100666666667
$ws.Cells.Item(8,"H").Value = 3.594302
$ws.Cells.Item(8,"I").Value = 0.2172419397410372
$ws.Cells.Item(8,"J").Value = 0.2211537917243501
$ws.Cells.Item(8,"K").Value = 1
$ws.Cells.Item(8,"L").Value = 0.3333333333333333
$ws.Cells.Item(8,"M").Value = 0.006995333333333333
$ws.Cells.Item(8,"N").Value = 0.020986
$ws.Cells.Item(8,"O").Value = 0.003071634989180631
$ws.Cells.Item(8,"P").Value = 0.003071634989180632
$ws.Cells.Item(8,"Q").Value = 0.008381113530222223
$ws.Cells.Item(8,"R").Value = 0.07543002177200001
$ws.Cells.Item(8,"S").Value = 0.0006672879432260402
$ws.Cells.Item(8,"T").Value = 0.0006793037246504798

# Row 9
$ws.Cells.Item(9,"A").Value = "Inflammatory-Mac"
$ws.Cells.Item(9,"B").Value = "Pglyrp1"
$ws.Cells.Item(9,"C").Value = "Trem1"
$ws.Cells.Item(9,"D").Value = "Inflammatory-Mac"
$ws.Cells.Item(9,"E").Value = 3
$ws.Cells.Item(9,"F").Value = 1
$ws.Cells.Item(9,"G").Value = 1.198100666666667
$ws.Cells.Item(9,"H").Value = 3.594302
$ws.Cells.Item(9,"I").Value = 0.2172419397410372
$ws.Cells.Item(9,"J").Value = 0.2211537917243501
$ws.Cells.Item(9,"K").Value = 3
$ws.Cells.Item(9,"L").Value = 1
$ws.Cells.Item(9,"M").Value = 1.994843666666667
$ws.Cells.Item(9,"N").Value = 5.984531
$ws.Cells.Item(9,"O").Value = 0.8759313262859122
$ws.Cells.Item(9,"P").Value = 0.8759313262859124
$ws.Cells.Item(9,"Q").Value = 2.390023526929112
$ws.Cells.Item(9,"R").Value = 21.510211742362
$ws.Cells.Item(9,"S").Value = 0.1902890204022909
$ws.Cells.Item(9,"T").Value = 0.1937155340982684

# Row 10
$ws.Cells.Item(10,"A").Value = "Inflammatory-Mac"
$ws.Cells.Item(10,"B").Value = "Pglyrp1"
$ws.Cells.Item(10,"C").Value = "Trem1"
$ws.Cells.Item(10,"D").Value = "Resolving-Mac"
$ws.Cells.Item(10,"E").Value = 3
$ws.Cells.Item(10,"F").Value = 1
$ws.Cells.Item(10,"G").Value = 1.198100666666667
$ws.Cells.Item(10,"H").Value = 3.594302
$ws.Cells.Item(10,"I").Value = 0.2172419397410372
$ws.Cells.Item(10,"J").Value = 0.2211537917243501
$ws.Cells.Item(10,"K").Value = 2
$ws.Cells.Item(10,"L").Value = 0.6666666666666666
$ws.Cells.Item(10,"M").Value = 0.2755583333333333
$ws.Cells.Item(10,"N").Value = 0.826675
$ws.Cells.Item(10,"O").Value = 0.120997038724907
$ws.Cells.Item(10,"P").Value = 0.120997038724907
$ws.Cells.Item(10,"Q").Value = 0.3301466228722222
$ws.Cells.Item(10,"R").Value = 2.97131960585
$ws.Cells.Item(10,"S").Value = 0.02628563139552019
$ws.Cells.Item(10,"T").Value = 0.02675895390143121

# Row 11
$ws.Cells.Item(11,"A").Value = "MuSCs"
$ws.Cells.Item(11,"B").Value = "Pglyrp1"
$ws.Cells.Item(11,"C").Value = "Trem1"
$ws.Cells.Item(11,"D").Value = "ECs"
$ws.Cells.Item(11,"E").Value = 2
$ws.Cells.Item(11,"F").Value = 1
$ws.Cells.Item(11,"G").Value = 0.292657
$ws.Cells.Item(11,"H").Value = 0.585314
$ws.Cells.Item(11,"I").Value = 0.05306513561642237
$ws.Cells.Item(11,"J").Value = 0.03601378249500076
$ws.Cells.Item(11,"K").Value = 1
$ws.Cells.Item(11,"L").Value = 0.3333333333333333
$ws.Cells.Item(11,"M").Value = 0.006995333333333333
$ws.Cells.Item(11,"N").Value = 0.020986
$ws.Cells.Item(11,"O").Value = 0.003071634989180631
$ws.Cells.Item(11,"P").Value = 0.003071634989180632
$ws.Cells.Item(11,"Q").Value = 0.002047233267333333
$ws.Cells.Item(11,"R").Value = 0.012283399604
$ws.Cells.Item(11,"S").Value = 0.0001629967272650183
$ws.Cells.Item(11,"T").Value = 0.0001106211944043853

# Row 12
$ws.Cells.Item(12,"A").Value = "MuSCs"
$ws.Cells.Item(12,"B").Value = "Pglyrp1"
$ws.Cells.Item(12,"C").Value = "Trem1"
$ws.Cells.Item(12,"D").Value = "Inflammatory-Mac"
$ws.Cells.Item(12,"E").Value = 2
$ws.Cells.Item(12,"F").Value = 1
$ws.Cells.Item(12,"G").Value = 0.292657
$ws.Cells.Item(12,"H").Value = 0.585314
$ws.Cells.Item(12,"I").Value = 0.05306513561642237
$ws.Cells.Item(12,"J").Value = 0.03601378249500076
$ws.Cells.Item(12,"K").Value = 3
$ws.Cells.Item(12,"L").Value = 1
$ws.Cells.Item(12,"M").Value = 1.994843666666667
$ws.Cells.Item(12,"N").Value = 5.984531
$ws.Cells.Item(12,"O").Value = 0.8759313262859122
$ws.Cells.Item(12,"P").Value = 0.8759313262859124
$ws.Cells.Item(12,"Q").Value = 0.5838049629556668
$ws.Cells.Item(12,"R").Value = 3.502829777734001
$ws.Cells.Item(12,"S").Value = 0.04648141462003465
$ws.Cells.Item(12,"T").Value = 0.0315456002654184

# Row 13
$ws.Cells.Item(13,"A").Value = "MuSCs"
$ws.Cells.Item(13,"B").Value = "Pglyrp1"
$ws.Cells.Item(13,"C").Value = "Trem1"
$ws.Cells.Item(13,"D").Value = "Resolving-Mac"
$ws.Cells.Item(13,"E").Value = 2
$ws.Cells.Item(13,"F").Value = 1
$ws.Cells.Item(13,"G").Value = 0.292657
$ws.Cells.Item(13,"H").Value = 0.585314
$ws.Cells.Item(13,"I").Value = 0.05306513561642237
$ws.Cells.Item(13,"J").Value = 0.03601378249500076
$ws.Cells.Item(13,"K").Value = 2
$ws.Cells.Item(13,"L").Value = 0.6666666666666666
$ws.Cells.Item(13,"M").Value = 0.2755583333333333
$ws.Cells.Item(13,"N").Value = 0.826675
$ws.Cells.Item(13,"O").Value = 0.120997038724907
$ws.Cells.Item(13,"P").Value = 0.120997038724907
$ws.Cells.Item(13,"Q").Value = 0.08064407515833334
$ws.Cells.Item(13,"R").Value = 0.48386445095
$ws.Cells.Item(13,"S").Value = 0.006420724269122699
$ws.Cells.Item(13,"T").Value = 0.004357561035177987

# Row 14
$ws.Cells.Item(14,"A").Value = "Resolving-Mac"
$ws.Cells.Item(14,"B").Value = "Pglyrp1"
$ws.Cells.Item(14,"C").Value = "Trem1"
$ws.Cells.Item(14,"D").Value = "ECs"
$ws.Cells.Item(14,"E").Value = 3
$ws.Cells.Item(14,"F").Value = 1
$ws.Cells.Item(14,"G").Value = 0.7183183333333333
$ws.Cells.Item(14,"H").Value = 2.154955
$ws.Cells.Item(14,"I").Value = 0.1302468752638612
$ws.Cells.Item(14,"J").Value = 0.1325922165820643
$ws.Cells.Item(14,"K").Value = 1
$ws.Cells.Item(14,"L").Value = 0.3333333333333333
$ws.Cells.Item(14,"M").Value = 0.006995333333333333
$ws.Cells.Item(14,"N").Value = 0.020986
$ws.Cells.Item(14,"O").Value = 0.003071634989180631
$ws.Cells.Item(14,"P").Value = 0.003071634989180632
$ws.Cells.Item(14,"Q").Value = 0.005024876181111111
$ws.Cells.Item(14,"R").Value = 0.04522388563
$ws.Cells.Item(14,"S").Value = 0.0004000708592919212
$ws.Cells.Item(14,"T").Value = 0.000407274891746485

# Row 15
$ws.Cells.Item(15,"A").Value = "Resolving-Mac"
$ws.Cells.Item(15,"B").Value = "Pglyrp1"
$ws.Cells.Item(15,"C").Value = "Trem1"
$ws.Cells.Item(15,"D").Value = "Inflammatory-Mac"
$ws.Cells.Item(15,"E").Value = 3
$ws.Cells.Item(15,"F").Value = 1
$ws.Cells.Item(15,"G").Value = 0.7183183333333333
$ws.Cells.Item(15,"H").Value = 2.154955
$ws.Cells.Item(15,"I").Value = 0.1302468752638612
$ws.Cells.Item(15,"J").Value = 0.1325922165820643
$ws.Cells.Item(15,"K").Value = 3
$ws.Cells.Item(15,"L").Value = 1
$ws.Cells.Item(15,"M").Value = 1.994843666666667
$ws.Cells.Item(15,"N").Value = 5.984531
$ws.Cells.Item(15,"O").Value = 0.8759313262859122
$ws.Cells.Item(15,"P").Value = 0.8759313262859124
$ws.Cells.Item(15,"Q").Value = 1.432932777900556
$ws.Cells.Item(15,"R").Value = 12.896395001105
$ws.Cells.Item(15,"S").Value = 0.1140873181944697
$ws.Cells.Item(15,"T").Value = 0.1161416761259165

# Row 16
$ws.Cells.Item(16,"A").Value = "Resolving-Mac"
$ws.Cells.Item(16,"B").Value = "Pglyrp1"
$ws.Cells.Item(16,"C").Value = "Trem1"
$ws.Cells.Item(16,"D").Value = "Resolving-Mac"
$ws.Cells.Item(16,"E").Value = 3
$ws.Cells.Item(16,"F").Value = 1
$ws.Cells.Item(16,"G").Value = 0.7183183333333333
$ws.Cells.Item(16,"H").Value = 2.154955
$ws.Cells.Item(16,"I").Value = 0.1302468752638612
$ws.Cells.Item(16,"J").Value = 0.1325922165820643
$ws.Cells.Item(16,"K").Value = 2
$ws.Cells.Item(16,"L").Value = 0.6666666666666666
$ws.Cells.Item(16,"M").Value = 0.2755583333333333
$ws.Cells.Item(16,"N").Value = 0.826675
$ws.Cells.Item(16,"O").Value = 0.120997038724907
$ws.Cells.Item(16,"P").Value = 0.120997038724907
$ws.Cells.Item(16,"Q").Value = 0.1979386027361111
$ws.Cells.Item(16,"R").Value = 1.781447424625
$ws.Cells.Item(16,"S").Value = 0.01575948621009954
$ws.Cells.Item(16,"T").Value = 0.01604326556440129
